$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.169.86'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '1.991.83'
$ws.Range('E3').Value = '  +5.78%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7754'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +63.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '254.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3465'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +19.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.70'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +23.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07031'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8428'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08177'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '100.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = '1.991.96'
$ws.Range('E14').Value = '  +5.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.615'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +15.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '272.19'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.35%  '
$ws.Range('D18').Value = '31.186.08'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008022'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.61%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.868'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.50%  '
$ws.Range('D21').Value = '2.252.65'
$ws.Range('E21').Value = '  +5.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.053'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.961'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1421'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +46.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.364'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +24.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.598'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.371'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.589'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.422'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.87%  '
$ws.Range('E34').Value = '  +8.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7881'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.214'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.758'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.904'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.685'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '79.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.61%  '
$ws.Range('E42').Value = '  +9.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.097'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8540'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.966'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.648'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.41%  '
$ws.Range('E50').Value = '  +16.04%  '
$ws.Range('E51').Value = '  +8.58%  '
